$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("M2").Value = 62.19572466666667
$ws.Range("N2").Value = 186.587174
$ws.Range("O2").Value = 0.9009169178676326
$ws.Range("P2").Value = 0.9009169178676325
$ws.Range("Q2").Value = 2937.526590858996
$ws.Range("R2").Value = 26437.73931773096
$ws.Range("S2").Value = 0.2922800945424654
$ws.Range("T2").Value = 0.2922800945424654
$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.03107099427955203
$ws.Range("P3").Value = 0.03107099427955203
$ws.Range("Q3").Value = 101.309976636516
$ws.Range("R3").Value = 911.7897897286439
$ws.Range("S3").Value = 0.01008021157716805
$ws.Range("T3").Value = 0.01008021157716805
$ws.Range("G4").Value = 47.23036199999999
$ws.Range("H4").Value = 141.691086
$ws.Range("I4").Value = 0.3244251370417807
$ws.Range("J4").Value = 0.3244251370417807
$ws.Range("O4").Value = 0.06801208785281536
$ws.Range("P4").Value = 0.06801208785281536
$ws.Range("Q4").Value = 221.759978756262
$ws.Range("R4").Value = 1995.839808806358
$ws.Range("S4").Value = 0.02206483092214725
$ws.Range("T4").Value = 0.02206483092214725
$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("M5").Value = 62.19572466666667
$ws.Range("N5").Value = 186.587174
$ws.Range("O5").Value = 0.9009169178676326
$ws.Range("P5").Value = 0.9009169178676325
$ws.Range("Q5").Value = 3792.54641434039
$ws.Range("R5").Value = 34132.91772906351
$ws.Range("S5").Value = 0.377353460557425
$ws.Range("T5").Value = 0.377353460557425
$ws.Range("I6").Value = 0.4188548944674916
$ws.Range("J6").Value = 0.4188548944674916
$ws.Range("O6").Value = 0.03107099427955203
$ws.Range("P6").Value = 0.03107099427955203
$ws.Range("S6").Value = 0.0130142380299618
$ws.Range("T6").Value = 0.0130142380299618
$ws.Range("I7").Value = 0.4188548944674916
$ws.Range("J7").Value = 0.4188548944674916
$ws.Range("O7").Value = 0.06801208785281536
$ws.Range("P7").Value = 0.06801208785281536
$ws.Range("S7").Value = 0.02848719588010474
$ws.Range("T7").Value = 0.02848719588010474
$ws.Range("I8").Value = 0.2567199684907278
$ws.Range("J8").Value = 0.2567199684907277
$ws.Range("M8").Value = 62.19572466666667
$ws.Range("N8").Value = 186.587174
$ws.Range("O8").Value = 0.9009169178676326
$ws.Range("P8").Value = 0.9009169178676325
$ws.Range("Q8").Value = 2324.486137918708
$ws.Range("R8").Value = 20920.37524126836
$ws.Range("S8").Value = 0.2312833627677422
$ws.Range("T8").Value = 0.2312833627677421
$ws.Range("I9").Value = 0.2567199684907278
$ws.Range("J9").Value = 0.2567199684907277
$ws.Range("O9").Value = 0.03107099427955203
$ws.Range("P9").Value = 0.03107099427955203
$ws.Range("S9").Value = 0.00797654467242218
$ws.Range("T9").Value = 0.007976544672422179
$ws.Range("O10").Value = 0.06801208785281536
$ws.Range("P10").Value = 0.06801208785281536
$ws.Range("S10").Value = 0.01746006105056337
$ws.Range("T10").Value = 0.01746006105056336
